$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Work from the BOTTOM of the document upward so that paragraph indices for
# not-yet-processed (earlier) paragraphs stay stable while we insert/delete.
# ---------------------------------------------------------------------------

# --- J. Delete the whole OC-3 section (paragraphs 24 .. 37 in the original
#        37-paragraph document): "OC-3: godkendVedligeholdelse" through
#        "gemVedligeholdelse er blevet vist til MidtTrafik".
$rStart = $d.Content
$rStart.Find.Execute("OC-3: godkendVedligeholdelse") | Out-Null
$startPos = $rStart.Start

$rEnd = $d.Content
$rEnd.Find.Execute("gemVedligeholdelse er blevet vist til MidtTrafik") | Out-Null
$endPos = $rEnd.End

$killRange = $d.Range($startPos, $endPos + 1)
$killRange.Delete()

# --- I. OC-2 "Slutbetingelser" hyperlink paragraph: change hyperlink text and
#        trailing text, then add a brand-new paragraph right after it.
$p23 = $d.Paragraphs(23)
$p23.Range.Find.Execute("ftp.bilTildelt", $true, $false, $false, $false, $false, $true, 1, $false, "ftp.angivKørselTilVEdligeholdelse", 2) | Out-Null

$p23 = $d.Paragraphs(23)
$p23.Range.Find.Execute("er blevet sat til at være sand", $true, $false, $false, $false, $false, $true, 1, $false, "er blevet kaldt.", 2) | Out-Null

$p23 = $d.Paragraphs(23)
$r = $p23.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$newP = $d.Paragraphs(24)
$newP.Range.InsertBefore("koerselsKartotek.gemGodkendtKoersel er blevet kaldt.")

# --- H. Delete "ftp.bil er blevet sat til bil" paragraph entirely.
$p22 = $d.Paragraphs(22)
$p22.Range.Delete()

# --- G. "En instans bil af Bil eksisterer" -> "En instans kørsel af Kørsel
#        eksisterer", then insert a brand-new paragraph right after it.
$p20 = $d.Paragraphs(20)
$p20.Range.Find.Execute("bil af Bil", $true, $false, $false, $false, $false, $true, 1, $false, "kørsel af Kørsel", 2) | Out-Null

$p20 = $d.Paragraphs(20)
$r = $p20.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$newP = $d.Paragraphs(21)
$newP.Range.InsertBefore("En instans koerselsKartotek af KoerselsKartotek eksisterer")

# --- F. OC-2 system-operation paragraph: "tildelBil" -> "angivKørselTilVedligeholdelse"
$p14 = $d.Paragraphs(14)
$p14.Range.Find.Execute("tildelBil", $true, $false, $false, $false, $false, $true, 1, $false, "angivKørselTilVedligeholdelse", 2) | Out-Null

# --- E. OC-2 heading: insert an empty Overskrift1 paragraph before it, then
#        change the heading text.
$p12 = $d.Paragraphs(12)
$p12.Range.InsertParagraphBefore()

$p13 = $d.Paragraphs(13)
$p13.Range.Find.Execute("tildelBil", $true, $false, $false, $false, $false, $true, 1, $false, "angivKørselTilVedligeholdelse", 2) | Out-Null

# --- D. OC-1 "Slutbetingelser" hyperlink paragraph: change hyperlink text and
#        trailing text.
$p11 = $d.Paragraphs(11)
$p11.Range.Find.Execute("ftp.kørsel", $true, $false, $false, $false, $false, $true, 1, $false, "ftp.tildelBil(bil)", 2) | Out-Null

$p11 = $d.Paragraphs(11)
$p11.Range.Find.Execute("er blevet sat til kørsel", $true, $false, $false, $false, $false, $true, 1, $false, "er blevet kaldt", 2) | Out-Null

# --- C. "En instans kørsel af Kørsel eksisterer" -> "En instans bil af Bil
#        eksisterer" (OC-1 section).
$p9 = $d.Paragraphs(9)
$p9.Range.Find.Execute("kørsel af Kørsel", $true, $false, $false, $false, $false, $true, 1, $false, "bil af Bil", 2) | Out-Null

# --- B. OC-1 system-operation paragraph: remove the _GoBack bookmark that
#        sits here (it gets re-created at the very end of the document) and
#        change the text.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$p3 = $d.Paragraphs(3)
$p3.Range.Find.Execute("angivKørselTilVedligeholdelse", $true, $false, $false, $false, $false, $true, 1, $false, "tildelBil", 2) | Out-Null

# --- A. OC-1 heading: insert an empty Overskrift1 paragraph before it, then
#        change the heading text.
$p1 = $d.Paragraphs(1)
$p1.Range.InsertParagraphBefore()

$p2 = $d.Paragraphs(2)
$p2.Range.Find.Execute("angivKørselTilVedligeholdelse", $true, $false, $false, $false, $false, $true, 1, $false, "tildelBil", 2) | Out-Null

# ---------------------------------------------------------------------------
# K. At the very end of the document, add a blank paragraph followed by a
#    paragraph that only carries the (re-created) _GoBack bookmark.
# ---------------------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

$endRange2 = $d.Content
$endRange2.Collapse(0)
$endRange2.InsertParagraphAfter()

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$bmRange = $lastPara.Range
$bmRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Host "done"
